# Auto-generated edits applying scheduled-runner price/profit updates
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 10
$ws.Range("I8").Value = 10
$ws.Range("K8").Value = 30
$ws.Range("M8").Value = 109

$ws.Range("H112").Value = 3226.0264
$ws.Range("J112").Value = 3385.6858
$ws.Range("L112").Value = 10157.0574
$ws.Range("N112").Value = -12373.0574

$ws.Range("H132").Value = 150890.66
$ws.Range("I132").Value = 349892.3
$ws.Range("J132").Value = 21539.574
$ws.Range("K132").Value = 1049676.9
$ws.Range("L132").Value = 64618.722
$ws.Range("M132").Value = -1047146.9
$ws.Range("N132").Value = -69678.72200000001

$ws.Range("H133").Value = 94500
$ws.Range("J133").Value = 94500
$ws.Range("L133").Value = 94500
$ws.Range("N133").Value = -104620

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 68584.3
$ws.Range("I45").Value = 82352.64
$ws.Range("J45").Value = 4332
$ws.Range("K45").Value = 82352.64
$ws.Range("L45").Value = 4332
$ws.Range("M45").Value = -81975.64
$ws.Range("N45").Value = -5086

$ws.Range("H61").Value = 10200.5
$ws.Range("I61").Value = 12257.857
$ws.Range("K61").Value = 12257.857
$ws.Range("M61").Value = -12045.857

$ws.Range("H102").Value = 654245.8
$ws.Range("I102").Value = 857798.3
$ws.Range("J102").Value = 2877.8
$ws.Range("K102").Value = 857798.3
$ws.Range("L102").Value = 2877.8
$ws.Range("M102").Value = -856176.3
$ws.Range("N102").Value = -6121.8

$ws.Range("H132").Value = 23141.133
$ws.Range("I132").Value = 28354
$ws.Range("J132").Value = 10977.777
$ws.Range("K132").Value = 85062
$ws.Range("L132").Value = 32933.331
$ws.Range("M132").Value = -82532
$ws.Range("N132").Value = -37993.331

$ws.Range("H136").Value = 10200.5
$ws.Range("I136").Value = 12257.857
$ws.Range("K136").Value = 36773.571
$ws.Range("M136").Value = -34223.571

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H36").Value = 1500
$ws.Range("I36").Value = 1500
$ws.Range("K36").Value = 1500
$ws.Range("M36").Value = -966

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()

$ws.Range("H7").Value = 362.125
$ws.Range("J7").Value = 378.14285
$ws.Range("L7").Value = 378.14285
$ws.Range("N7").Value = -604.14285

$ws.Range("H58").Value = 459124.1
$ws.Range("I58").Value = 1251431.2
$ws.Range("J58").Value = 6377.143
$ws.Range("K58").Value = 1251431.2
$ws.Range("L58").Value = 6377.143
$ws.Range("M58").Value = -1251228.2
$ws.Range("N58").Value = -6783.143

$ws.Range("H99").Value = 10614.828
$ws.Range("I99").Value = 11174.565
$ws.Range("K99").Value = 11174.565
$ws.Range("M99").Value = -9676.565000000001

$ws.Range("H119").Value = 75000
$ws.Range("J119").Value = 75000
$ws.Range("L119").Value = 75000
$ws.Range("N119").Value = -84676

$ws.Range("H122").Value = 4062.15
$ws.Range("I122").Value = 2537.923
$ws.Range("J122").Value = 6892.857
$ws.Range("K122").Value = 7613.768999999999
$ws.Range("L122").Value = 20678.571
$ws.Range("M122").Value = -5163.768999999999
$ws.Range("N122").Value = -25578.571

$ws.Range("H126").Value = 10614.828
$ws.Range("I126").Value = 11174.565
$ws.Range("K126").Value = 33523.695
$ws.Range("M126").Value = -31053.695

$ws.Range("H132").Value = 15162757
$ws.Range("I132").Value = 17556526
$ws.Range("K132").Value = 52669578
$ws.Range("M132").Value = -52667048

$ws.Range("H136").Value = 459124.1
$ws.Range("I136").Value = 1251431.2
$ws.Range("J136").Value = 6377.143
$ws.Range("K136").Value = 3754293.6
$ws.Range("L136").Value = 19131.429
$ws.Range("M136").Value = -3751743.6
$ws.Range("N136").Value = -24231.429

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 14887236
$ws.Range("I131").Value = 13890942
$ws.Range("J131").Value = 15158952
$ws.Range("K131").Value = 41672826
$ws.Range("L131").Value = 45476856
$ws.Range("M131").Value = -41667786
$ws.Range("N131").Value = -45486936

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 19239048
$ws.Range("J102").Value = 8149.75
$ws.Range("L102").Value = 8149.75
$ws.Range("N102").Value = -11393.75

$ws.Range("H122").Value = 412257.78
$ws.Range("I122").Value = 553198.2
$ws.Range("J122").Value = 9570.857
$ws.Range("K122").Value = 1659594.6
$ws.Range("L122").Value = 28712.571
$ws.Range("M122").Value = -1657144.6
$ws.Range("N122").Value = -33612.571

$ws.Range("H126").Value = 4745.7393
$ws.Range("I126").Value = 2815.9167
$ws.Range("K126").Value = 8447.750100000001
$ws.Range("M126").Value = -5977.750100000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6833.8335
$ws.Range("I7").Value = 6000.5
$ws.Range("J7").Value = 7250.5
$ws.Range("K7").Value = 6000.5
$ws.Range("L7").Value = 7250.5
$ws.Range("M7").Value = -5888.5
$ws.Range("N7").Value = -7474.5

$ws.Range("H40").Value = 66669290
$ws.Range("I40").Value = 3034
$ws.Range("J40").Value = 166668670
$ws.Range("K40").Value = 3034
$ws.Range("L40").Value = 166668670
$ws.Range("M40").Value = -2898
$ws.Range("N40").Value = -166668942

$ws.Range("H46").Value = 6139.2334
$ws.Range("J46").Value = 6220.607
$ws.Range("L46").Value = 6220.607
$ws.Range("N46").Value = -6596.607

$ws.Range("H100").Value = 4893
$ws.Range("I100").Value = 4893
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 4893
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -4352
$ws.Range("N100").ClearContents()

$ws.Range("H126").Value = 6833.8335
$ws.Range("I126").Value = 6000.5
$ws.Range("J126").Value = 7250.5
$ws.Range("K126").Value = 18001.5
$ws.Range("L126").Value = 21751.5
$ws.Range("M126").Value = -15531.5
$ws.Range("N126").Value = -26691.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2920.3333
$ws.Range("I126").Value = 2555.6
$ws.Range("J126").Value = 3649.8
$ws.Range("K126").Value = 7666.799999999999
$ws.Range("L126").Value = 10949.4
$ws.Range("M126").Value = -5196.799999999999
$ws.Range("N126").Value = -15889.4

$ws.Range("H136").Value = 8654.369000000001
$ws.Range("I136").Value = 2433.9092
$ws.Range("J136").Value = 11836.931
$ws.Range("K136").Value = 7301.7276
$ws.Range("L136").Value = 35510.79300000001
$ws.Range("M136").Value = -4751.7276
$ws.Range("N136").Value = -40610.79300000001
